# Autogenerated on Fri Mar 20 2015 00:16:06 GMT+0000 (Coordinated Universal Time)
# Remove the "Enterprises density (per 1000 people)" row from the MSME
# summary sheet and let everything below it shift up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("14").Delete()
